$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.894.74"
$ws.Range("E2").Value = "  +1.20%  "
$ws.Range("D3").Value = "3.516.80"
$ws.Range("E3").Value = "  -0.02%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "605.68"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +4.30%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "171.63"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.78%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.616"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -1.35%  "
$ws.Range("D8").Value = "3.513.22"
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("E9").Value = "  -0.12%  "
$ws.Range("E10").Value = "  +4.67%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.67"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.76%  "
$ws.Range("E12").Value = "  -3.03%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "47.25"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.22%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000280"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.81%  "
$ws.Range("D15").Value = "4.080.27"
$ws.Range("E15").Value = "  -0.18%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "620.31"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -8.20%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "8.39"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -4.59%  "
$ws.Range("D18").Value = "3.515.53"
$ws.Range("E18").Value = "  -0.37%  "
$ws.Range("D19").Value = "69.886.86"
$ws.Range("E19").Value = "  +1.18%  "
$ws.Range("E20").Value = "  -2.04%  "
$ws.Range("E21").Value = "  -1.57%  "
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.01"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -11.17%  "
$ws.Range("B23").Value = "Polygon"
$ws.Range("C23").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.885"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -2.65%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "15.76"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -2.85%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "96.09"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -2.23%  "
$ws.Range("E26").Value = "  -0.60%  "
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("E28").Value = "  -2.89%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.24"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -2.73%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "33.16"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.39%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.44"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -3.90%  "
$ws.Range("E32").Value = "  -4.44%  "
$ws.Range("E33").Value = "  -1.59%  "
$ws.Range("E34").Value = "  -5.90%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "569.24"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -1.91%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.77"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -1.42%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.49"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -3.31%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "56.99"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.70%  "
$ws.Range("E39").Value = "  -3.85%  "
$ws.Range("E40").Value = "  +0.04%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0451"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +2.13%  "
$ws.Range("E42").Value = "  +2.71%  "
$ws.Range("E43").Value = "  -3.70%  "
$ws.Range("D44").Value = "3.333.70"
$ws.Range("E44").Value = "  -2.97%  "
$ws.Range("B45").Value = "ThetaToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.98"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +2.31%  "
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "33.13"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -1.42%  "
$ws.Range("D47").Value = "0.0₃0704"
$ws.Range("E47").Value = "  -0.81%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.63"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.71%  "
$ws.Range("E49").Value = "  -3.60%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "136.12"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +3.18%  "
$ws.Range("E51").Value = "  +0.97%  "
